# Apply "Added ability to read in Light Data from blender" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the earned value for "Support for adding Point and/or Spotlight
# sources via Blender" (row 15) - scored 0.07 out of 0.11.
$ws.Range("C15").Value = 0.07
# Match the "filled in / graded" look used by the other completed rows
# (red font percentage style, same as C14, C47, etc.)
$ws.Range("C15").Font.Color = $ws.Range("C14").Font.Color
$ws.Range("C15").NumberFormat = $ws.Range("C14").NumberFormat

# Fill in the earned value for "Manage your project in a private GIT repo
# (GitHub or GitLab)" (row 46) - scored full 0.02.
$ws.Range("C46").Value = 0.02

# Update the window view / selection state to reflect where the user was
# working when they made the edit (scrolled down so row 10 is at the top,
# with A15 selected).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A15").Select() | Out-Null

$excel.ActiveWindow.Left = 3420
$excel.ActiveWindow.Top = 1725
$excel.ActiveWindow.Width = 21600
$excel.ActiveWindow.Height = 11295

$wb.Save()
